$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3447.5
$ws.Range("I64").Value = 3566.6667
$ws.Range("J64").Value = 3090
$ws.Range("K64").Value = 3566.6667
$ws.Range("L64").Value = 3090
$ws.Range("M64").Value = -3318.6667
$ws.Range("N64").Value = -3586

$ws.Range("H67").Value = 3447.5
$ws.Range("I67").Value = 3566.6667
$ws.Range("J67").Value = 3090
$ws.Range("K67").Value = 3566.6667
$ws.Range("L67").Value = 3090
$ws.Range("M67").Value = -2708.6667
$ws.Range("N67").Value = -4806

$ws.Range("H98").Value = 7426.125
$ws.Range("I98").Value = 9068.166999999999
$ws.Range("K98").Value = 9068.166999999999
$ws.Range("M98").Value = -7570.166999999999

$ws.Range("H116").Value = 4175.9
$ws.Range("I116").Value = 3618
$ws.Range("J116").Value = 4857.778
$ws.Range("K116").Value = 3618
$ws.Range("L116").Value = 4857.778
$ws.Range("M116").Value = -176
$ws.Range("N116").Value = -11741.778

$ws.Range("H122").Value = 7426.125
$ws.Range("I122").Value = 9068.166999999999
$ws.Range("K122").Value = 27204.501
$ws.Range("M122").Value = -24754.501

$ws.Range("H132").Value = 3638224.5
$ws.Range("I132").Value = 4083366.2
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 12250098.6
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -12247568.6
$ws.Range("N132").Value = -13760

$ws.Range("H133").Value = 16339.143
$ws.Range("J133").Value = 16339.143
$ws.Range("L133").Value = 16339.143
$ws.Range("N133").Value = -26459.143

$ws.Range("H135").Value = 1414.1333
$ws.Range("I135").Value = 1542.4546
$ws.Range("J135").Value = 1061.25
$ws.Range("K135").Value = 13882.0914
$ws.Range("L135").Value = 9551.25
$ws.Range("M135").Value = -11347.0914
$ws.Range("N135").Value = -14621.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5042
$ws.Range("I61").Value = 8096.25
$ws.Range("J61").Value = 2167.4119
$ws.Range("K61").Value = 8096.25
$ws.Range("L61").Value = 2167.4119
$ws.Range("M61").Value = -7884.25
$ws.Range("N61").Value = -2591.4119

$ws.Range("H74").Value = 2141.9092
$ws.Range("I74").Value = 1392.2142
$ws.Range("J74").Value = 3453.875
$ws.Range("K74").Value = 1392.2142
$ws.Range("L74").Value = 3453.875
$ws.Range("M74").Value = -518.2141999999999
$ws.Range("N74").Value = -5201.875

$ws.Range("H77").Value = 2141.9092
$ws.Range("I77").Value = 1392.2142
$ws.Range("J77").Value = 3453.875
$ws.Range("K77").Value = 6961.071
$ws.Range("L77").Value = 17269.375
$ws.Range("M77").Value = -2593.071
$ws.Range("N77").Value = -26005.375

$ws.Range("H97").Value = 1227.1818
$ws.Range("I97").Value = 1055.4445
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1055.4445
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -559.4445000000001
$ws.Range("N97").Value = -2992

$ws.Range("H122").Value = 1852.2
$ws.Range("I122").Value = 1612.4584
$ws.Range("J122").Value = 2375.2727
$ws.Range("K122").Value = 4837.3752
$ws.Range("L122").Value = 7125.8181
$ws.Range("M122").Value = -2387.3752
$ws.Range("N122").Value = -12025.8181

$ws.Range("H132").Value = 3159.3171
$ws.Range("I132").Value = 3197.3215
$ws.Range("J132").Value = 3077.4614
$ws.Range("K132").Value = 9591.9645
$ws.Range("L132").Value = 9232.3842
$ws.Range("M132").Value = -7061.9645
$ws.Range("N132").Value = -14292.3842

$ws.Range("H136").Value = 5042
$ws.Range("I136").Value = 8096.25
$ws.Range("J136").Value = 2167.4119
$ws.Range("K136").Value = 24288.75
$ws.Range("L136").Value = 6502.2357
$ws.Range("M136").Value = -21738.75
$ws.Range("N136").Value = -11602.2357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 557.1429000000001
$ws.Range("I22").Value = 520
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 520
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -347
$ws.Range("N22").Value = -996

$ws.Range("H86").Value = 87796.664
$ws.Range("I86").Value = 7440
$ws.Range("J86").Value = 127975
$ws.Range("K86").Value = 7440
$ws.Range("L86").Value = 127975
$ws.Range("M86").Value = -6317
$ws.Range("N86").Value = -130221

$ws.Range("H89").Value = 87796.664
$ws.Range("I89").Value = 7440
$ws.Range("J89").Value = 127975
$ws.Range("K89").Value = 37200
$ws.Range("L89").Value = 639875
$ws.Range("M89").Value = -31584
$ws.Range("N89").Value = -651107

$ws.Range("H94").Value = 515.6
$ws.Range("I94").Value = 352.6
$ws.Range("K94").Value = 352.6
$ws.Range("M94").Value = 98.39999999999998

$ws.Range("H134").Value = 2314.6843
$ws.Range("I134").Value = 2471.8125
$ws.Range("J134").Value = 1476.6666
$ws.Range("K134").Value = 7415.4375
$ws.Range("L134").Value = 4429.9998
$ws.Range("M134").Value = -4880.4375
$ws.Range("N134").Value = -9499.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2372.6086
$ws.Range("I31").Value = 1743.6
$ws.Range("J31").Value = 4374
$ws.Range("K31").Value = 1743.6
$ws.Range("L31").Value = 4374
$ws.Range("M31").Value = -1448.6
$ws.Range("N31").Value = -4964

$ws.Range("H34").Value = 2372.6086
$ws.Range("I34").Value = 1743.6
$ws.Range("J34").Value = 4374
$ws.Range("K34").Value = 1743.6
$ws.Range("L34").Value = 4374
$ws.Range("M34").Value = -1541.6
$ws.Range("N34").Value = -4778

$ws.Range("H94").Value = 3471.2
$ws.Range("I94").Value = 9460.75
$ws.Range("J94").Value = 2330.3333
$ws.Range("K94").Value = 9460.75
$ws.Range("L94").Value = 2330.3333
$ws.Range("M94").Value = -9009.75
$ws.Range("N94").Value = -3232.3333

$ws.Range("H107").Value = 599.6774
$ws.Range("I107").Value = 684.2857
$ws.Range("J107").Value = 422
$ws.Range("K107").Value = 684.2857
$ws.Range("L107").Value = 422
$ws.Range("M107").Value = 1235.7143
$ws.Range("N107").Value = -4262

$ws.Range("H132").Value = 1901.1875
$ws.Range("I132").Value = 2078.8635
$ws.Range("J132").Value = 1510.3
$ws.Range("K132").Value = 6236.5905
$ws.Range("L132").Value = 4530.9
$ws.Range("M132").Value = -3706.5905
$ws.Range("N132").Value = -9590.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1195.9445
$ws.Range("J131").Value = 1009.5085
$ws.Range("L131").Value = 3028.5255
$ws.Range("N131").Value = -13108.5255

$ws.Range("H132").Value = 1520.8334
$ws.Range("I132").Value = 1531.25
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 13781.25
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -11251.25
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 28500
$ws.Range("J69").Value = 28500
$ws.Range("L69").Value = 28500
$ws.Range("N69").Value = -29998

$ws.Range("H72").Value = 28500
$ws.Range("J72").Value = 28500
$ws.Range("L72").Value = 85500
$ws.Range("N72").Value = -92988

$ws.Range("H80").Value = 4266.6665
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 4266.6665
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

$ws.Range("H97").Value = 1571.3334
$ws.Range("I97").Value = 1590.65
$ws.Range("J97").Value = 1474.75
$ws.Range("K97").Value = 1590.65
$ws.Range("L97").Value = 1474.75
$ws.Range("M97").Value = -1094.65
$ws.Range("N97").Value = -2466.75

$ws.Range("H122").Value = 6950.8
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 7786.857
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 23360.571
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -28260.571

$ws.Range("H132").Value = 3432.4634
$ws.Range("I132").Value = 3880
$ws.Range("J132").Value = 3045.9546
$ws.Range("K132").Value = 11640
$ws.Range("L132").Value = 9137.863799999999
$ws.Range("M132").Value = -9110
$ws.Range("N132").Value = -14197.8638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1869.1818
$ws.Range("I40").Value = 1760.1666
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1760.1666
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1624.1666
$ws.Range("N40").Value = -2272

$ws.Range("H61").Value = 1633.3529
$ws.Range("I61").Value = 519.0714
$ws.Range("J61").Value = 6833.3335
$ws.Range("K61").Value = 519.0714
$ws.Range("L61").Value = 6833.3335
$ws.Range("M61").Value = -317.0714
$ws.Range("N61").Value = -7237.3335

$ws.Range("H97").Value = 20274.777
$ws.Range("J97").Value = 20274.777
$ws.Range("L97").Value = 20274.777
$ws.Range("N97").Value = -22256.777

$ws.Range("H113").Value = 1633.3529
$ws.Range("I113").Value = 519.0714
$ws.Range("J113").Value = 6833.3335
$ws.Range("K113").Value = 519.0714
$ws.Range("L113").Value = 6833.3335
$ws.Range("M113").Value = 1650.9286
$ws.Range("N113").Value = -11173.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 796.0417
$ws.Range("I113").Value = 417.8889
$ws.Range("J113").Value = 1930.5
$ws.Range("K113").Value = 1253.6667
$ws.Range("L113").Value = 5791.5
$ws.Range("M113").Value = 916.3333
$ws.Range("N113").Value = -10131.5

$ws.Range("H132").Value = 4083.3062
$ws.Range("I132").Value = 1716.561
$ws.Range("J132").Value = 16212.875
$ws.Range("K132").Value = 5149.683
$ws.Range("L132").Value = 48638.625
$ws.Range("M132").Value = -2619.683
$ws.Range("N132").Value = -53698.625

$ws.Range("H136").Value = 3532.5625
$ws.Range("I136").Value = 3792.4
$ws.Range("J136").Value = 3099.5
$ws.Range("K136").Value = 11377.2
$ws.Range("L136").Value = 9298.5
$ws.Range("M136").Value = -8827.200000000001
$ws.Range("N136").Value = -14398.5
